# Insert a new translation row ("strFileHeader29") into the Slovak
# (sk-SK) localization table, just above the existing "strFileHeaderSection"
# row (which currently lives at row 50), shifting all following rows down
# by one.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a blank worksheet row at row 50, pushing everything from the old
# row 50 onward down to row 51 onward.
$ws.Rows("50:50").Insert()

# The table (ListObject) range does not auto-grow when rows are inserted
# this way, so resize it to include the new row and keep the total row
# count consistent with the new last data row (159).
$tbl.Resize($ws.Range("B2:E159"))

# Populate the new row's cells with the new key/comment/value strings.
$ws.Range("B50").Value = "strFileHeader29"
$ws.Range("C50").Value = "Field description in exported file"
$ws.Range("D50").Value = "Differentiation algorithm"

# Key column (B): left/top aligned, no wrap text.
$keyCell = $ws.Range("B50")
$keyCell.HorizontalAlignment = -4131
$keyCell.VerticalAlignment = -4108
$keyCell.WrapText = $false

# Comment/value columns (C:D): same alignment, no wrap text, and unlocked
# (so translators can edit them once the sheet is protected).
$dataCells = $ws.Range("C50:D50")
$dataCells.HorizontalAlignment = -4131
$dataCells.VerticalAlignment = -4108
$dataCells.WrapText = $false
$dataCells.Locked = $false
